# Restores the accelerometer sample sheet to A1:C31, pushing the previous
# 20 rows of samples down and filling rows 2-17 with the new May 9th capture.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 30,3
$data[0,0] = -1.181560516357422; $data[0,1] = 1.240549564361572; $data[0,2] = -0.5883253216743469
$data[1,0] = -1.229650020599365; $data[1,1] = 1.595798969268799; $data[1,2] = -0.7779194116592407
$data[2,0] = -0.2790107727050781; $data[2,1] = 1.227867126464844; $data[2,2] = -1.016065120697022
$data[3,0] = -1.14784574508667; $data[3,1] = 1.330484867095947; $data[3,2] = -0.8695551156997681
$data[4,0] = -1.288122177124023; $data[4,1] = 1.133394718170166; $data[4,2] = -0.7701815962791443
$data[5,0] = -1.425555229187011; $data[5,1] = 1.231297016143799; $data[5,2] = -0.6422767043113708
$data[6,0] = -1.529905319213867; $data[6,1] = 1.149116039276123; $data[6,2] = -0.9077204465866088
$data[7,0] = -1.152251720428467; $data[7,1] = 1.056320667266846; $data[7,2] = -0.8245069980621338
$data[8,0] = -1.123115062713623; $data[8,1] = 1.085736274719239; $data[8,2] = -0.7493376731872559
$data[9,0] = -1.394256114959717; $data[9,1] = 1.270269870758057; $data[9,2] = -0.8571128845214844
$data[10,0] = -1.345842838287354; $data[10,1] = 1.367433071136475; $data[10,2] = -0.7304041385650635
$data[11,0] = -1.356554985046387; $data[11,1] = 1.367568492889404; $data[11,2] = -0.6971900463104248
$data[12,0] = -1.285661697387695; $data[12,1] = 1.301629066467285; $data[12,2] = -0.5075737237930298
$data[13,0] = -1.321344375610352; $data[13,1] = 1.155209541320801; $data[13,2] = -0.5763433575630188
$data[14,0] = -1.26480770111084; $data[14,1] = 1.153616428375244; $data[14,2] = -0.6878960132598877
$data[15,0] = -0.8195595741271973; $data[15,1] = 1.118548393249511; $data[15,2] = -0.9026113748550416
$data[16,0] = -1.078746795654297; $data[16,1] = 1.204861640930176; $data[16,2] = -0.6440812349319458
$data[17,0] = -1.571155071258545; $data[17,1] = 1.288459777832031; $data[17,2] = -0.5561246871948242
$data[18,0] = -1.56222677230835; $data[18,1] = 1.121483325958252; $data[18,2] = -0.1114475727081298
$data[19,0] = -1.141444206237793; $data[19,1] = 1.122160911560059; $data[19,2] = 0.2474624365568161
$data[20,0] = -1.447634696960449; $data[20,1] = 0.8090605735778809; $data[20,2] = 2.089949369430542
$data[21,0] = -2.48740816116333; $data[21,1] = -0.6178178787231445; $data[21,2] = 2.327399969100952
$data[22,0] = -2.30517578125; $data[22,1] = -0.4891290664672851; $data[22,2] = 2.537726402282715
$data[23,0] = -3.007162570953369; $data[23,1] = -1.502803325653076; $data[23,2] = 3.441044807434082
$data[24,0] = -3.378014087677002; $data[24,1] = -3.04712963104248; $data[24,2] = 5.235954284667969
$data[25,0] = -6.576094627380371; $data[25,1] = 0.4697372913360595; $data[25,2] = 2.406523942947388
$data[26,0] = 8.52888298034668; $data[26,1] = 7.144853115081787; $data[26,2] = -10.8126277923584
$data[27,0] = -10.79609775543213; $data[27,1] = -1.548766374588013; $data[27,2] = -0.7228314876556396
$data[28,0] = -0.0029077529907226; $data[28,1] = 0.2878659963607788; $data[28,2] = 0.002394676208496
$data[29,0] = -3.11443567276001; $data[29,1] = 0.3510211706161499; $data[29,2] = -0.4145855903625488

$ws.Range("A2:C31").Value = $data
